$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the ANTIOQUIA / MEDELLIN row (row 2); CAUCA/TIMBIQUI shifts up to row 2
$ws.Rows.Item(2).Delete()
